$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Set widths for new columns H:K to match column G as closely as possible
$ws.Columns.Item(8).ColumnWidth = 8.67
$ws.Columns.Item(9).ColumnWidth = 8.67
$ws.Columns.Item(10).ColumnWidth = 8.67
$ws.Columns.Item(11).ColumnWidth = 8.67

# Block starting at row 1
$ws.Range("G1").Copy()
$ws.Range("H1:K1").PasteSpecial(-4122)
$ws.Range("G1").Value = 2024
$ws.Range("H1").Value = 2025
$ws.Range("I1").Value = 2026
$ws.Range("J1").Value = 2027
$ws.Range("K1").Value = 2028
$ws.Range("G2").Copy()
$ws.Range("H2:K6").PasteSpecial(-4122)

# Block starting at row 8
$ws.Range("G8").Copy()
$ws.Range("H8:K8").PasteSpecial(-4122)
$ws.Range("G8").Value = 2024
$ws.Range("H8").Value = 2025
$ws.Range("I8").Value = 2026
$ws.Range("J8").Value = 2027
$ws.Range("K8").Value = 2028
$ws.Range("G9").Copy()
$ws.Range("H9:K13").PasteSpecial(-4122)

# Block starting at row 15
$ws.Range("G15").Copy()
$ws.Range("H15:K15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2024
$ws.Range("H15").Value = 2025
$ws.Range("I15").Value = 2026
$ws.Range("J15").Value = 2027
$ws.Range("K15").Value = 2028
$ws.Range("G16").Copy()
$ws.Range("H16:K20").PasteSpecial(-4122)

# Block starting at row 22
$ws.Range("G22").Copy()
$ws.Range("H22:K22").PasteSpecial(-4122)
$ws.Range("G22").Value = 2024
$ws.Range("H22").Value = 2025
$ws.Range("I22").Value = 2026
$ws.Range("J22").Value = 2027
$ws.Range("K22").Value = 2028
$ws.Range("G23").Copy()
$ws.Range("H23:K27").PasteSpecial(-4122)

# Block starting at row 29
$ws.Range("G29").Copy()
$ws.Range("H29:K29").PasteSpecial(-4122)
$ws.Range("G29").Value = 2024
$ws.Range("H29").Value = 2025
$ws.Range("I29").Value = 2026
$ws.Range("J29").Value = 2027
$ws.Range("K29").Value = 2028
$ws.Range("G30").Copy()
$ws.Range("H30:K34").PasteSpecial(-4122)

# Block starting at row 36
$ws.Range("G36").Copy()
$ws.Range("H36:K36").PasteSpecial(-4122)
$ws.Range("G36").Value = 2024
$ws.Range("H36").Value = 2025
$ws.Range("I36").Value = 2026
$ws.Range("J36").Value = 2027
$ws.Range("K36").Value = 2028
$ws.Range("G37").Copy()
$ws.Range("H37:K41").PasteSpecial(-4122)

# Block starting at row 43
$ws.Range("G43").Copy()
$ws.Range("H43:K43").PasteSpecial(-4122)
$ws.Range("G43").Value = 2024
$ws.Range("H43").Value = 2025
$ws.Range("I43").Value = 2026
$ws.Range("J43").Value = 2027
$ws.Range("K43").Value = 2028
$ws.Range("G44").Copy()
$ws.Range("H44:K48").PasteSpecial(-4122)

# Block starting at row 50
$ws.Range("G50").Copy()
$ws.Range("H50:K50").PasteSpecial(-4122)
$ws.Range("G50").Value = 2024
$ws.Range("H50").Value = 2025
$ws.Range("I50").Value = 2026
$ws.Range("J50").Value = 2027
$ws.Range("K50").Value = 2028
$ws.Range("G51").Copy()
$ws.Range("H51:K55").PasteSpecial(-4122)

# Block starting at row 57
$ws.Range("G57").Copy()
$ws.Range("H57:K57").PasteSpecial(-4122)
$ws.Range("G57").Value = 2024
$ws.Range("H57").Value = 2025
$ws.Range("I57").Value = 2026
$ws.Range("J57").Value = 2027
$ws.Range("K57").Value = 2028
$ws.Range("G58").Copy()
$ws.Range("H58:K62").PasteSpecial(-4122)

# Block starting at row 64
$ws.Range("G64").Copy()
$ws.Range("H64:K64").PasteSpecial(-4122)
$ws.Range("G64").Value = 2024
$ws.Range("H64").Value = 2025
$ws.Range("I64").Value = 2026
$ws.Range("J64").Value = 2027
$ws.Range("K64").Value = 2028
$ws.Range("G65").Copy()
$ws.Range("H65:K69").PasteSpecial(-4122)

# Block starting at row 71
$ws.Range("G71").Copy()
$ws.Range("H71:K71").PasteSpecial(-4122)
$ws.Range("G71").Value = 2024
$ws.Range("H71").Value = 2025
$ws.Range("I71").Value = 2026
$ws.Range("J71").Value = 2027
$ws.Range("K71").Value = 2028
$ws.Range("G72").Copy()
$ws.Range("H72:K76").PasteSpecial(-4122)

# Block starting at row 78
$ws.Range("G78").Copy()
$ws.Range("H78:K78").PasteSpecial(-4122)
$ws.Range("G78").Value = 2024
$ws.Range("H78").Value = 2025
$ws.Range("I78").Value = 2026
$ws.Range("J78").Value = 2027
$ws.Range("K78").Value = 2028
$ws.Range("G79").Copy()
$ws.Range("H79:K83").PasteSpecial(-4122)

# Block starting at row 85
$ws.Range("G85").Copy()
$ws.Range("H85:K85").PasteSpecial(-4122)
$ws.Range("G85").Value = 2024
$ws.Range("H85").Value = 2025
$ws.Range("I85").Value = 2026
$ws.Range("J85").Value = 2027
$ws.Range("K85").Value = 2028
$ws.Range("G86").Copy()
$ws.Range("H86:K90").PasteSpecial(-4122)

# Update conditional formatting formulas from G{row}:G{row} to G{row}:K{row}
$cf = $ws.Range("E2").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G2:K2,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G2:K2,"<>" & "")>0,NOT(ISBLANK(E2)))'
$cf = $ws.Range("E3").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G3:K3,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G3:K3,"<>" & "")>0,NOT(ISBLANK(E3)))'
$cf = $ws.Range("E4").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G4:K4,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G4:K4,"<>" & "")>0,NOT(ISBLANK(E4)))'
$cf = $ws.Range("E5").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G5:K5,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G5:K5,"<>" & "")>0,NOT(ISBLANK(E5)))'
$cf = $ws.Range("E6").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G6:K6,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G6:K6,"<>" & "")>0,NOT(ISBLANK(E6)))'
$cf = $ws.Range("E9").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G9:K9,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G9:K9,"<>" & "")>0,NOT(ISBLANK(E9)))'
$cf = $ws.Range("E10").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G10:K10,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G10:K10,"<>" & "")>0,NOT(ISBLANK(E10)))'
$cf = $ws.Range("E11").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G11:K11,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G11:K11,"<>" & "")>0,NOT(ISBLANK(E11)))'
$cf = $ws.Range("E12").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G12:K12,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G12:K12,"<>" & "")>0,NOT(ISBLANK(E12)))'
$cf = $ws.Range("E13").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G13:K13,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G13:K13,"<>" & "")>0,NOT(ISBLANK(E13)))'
$cf = $ws.Range("E16").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G16:K16,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G16:K16,"<>" & "")>0,NOT(ISBLANK(E16)))'
$cf = $ws.Range("E17").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G17:K17,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G17:K17,"<>" & "")>0,NOT(ISBLANK(E17)))'
$cf = $ws.Range("E18").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G18:K18,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G18:K18,"<>" & "")>0,NOT(ISBLANK(E18)))'
$cf = $ws.Range("E19").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G19:K19,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G19:K19,"<>" & "")>0,NOT(ISBLANK(E19)))'
$cf = $ws.Range("E20").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G20:K20,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G20:K20,"<>" & "")>0,NOT(ISBLANK(E20)))'
$cf = $ws.Range("E23").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G23:K23,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G23:K23,"<>" & "")>0,NOT(ISBLANK(E23)))'
$cf = $ws.Range("E24").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G24:K24,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G24:K24,"<>" & "")>0,NOT(ISBLANK(E24)))'
$cf = $ws.Range("E25").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G25:K25,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G25:K25,"<>" & "")>0,NOT(ISBLANK(E25)))'
$cf = $ws.Range("E26").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G26:K26,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G26:K26,"<>" & "")>0,NOT(ISBLANK(E26)))'
$cf = $ws.Range("E27").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G27:K27,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G27:K27,"<>" & "")>0,NOT(ISBLANK(E27)))'
$cf = $ws.Range("E30").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G30:K30,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G30:K30,"<>" & "")>0,NOT(ISBLANK(E30)))'
$cf = $ws.Range("E31").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G31:K31,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G31:K31,"<>" & "")>0,NOT(ISBLANK(E31)))'
$cf = $ws.Range("E32").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G32:K32,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G32:K32,"<>" & "")>0,NOT(ISBLANK(E32)))'
$cf = $ws.Range("E33").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G33:K33,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G33:K33,"<>" & "")>0,NOT(ISBLANK(E33)))'
$cf = $ws.Range("E34").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G34:K34,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G34:K34,"<>" & "")>0,NOT(ISBLANK(E34)))'
$cf = $ws.Range("E37").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G37:K37,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G37:K37,"<>" & "")>0,NOT(ISBLANK(E37)))'
$cf = $ws.Range("E38").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G38:K38,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G38:K38,"<>" & "")>0,NOT(ISBLANK(E38)))'
$cf = $ws.Range("E39").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G39:K39,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G39:K39,"<>" & "")>0,NOT(ISBLANK(E39)))'
$cf = $ws.Range("E40").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G40:K40,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G40:K40,"<>" & "")>0,NOT(ISBLANK(E40)))'
$cf = $ws.Range("E41").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G41:K41,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G41:K41,"<>" & "")>0,NOT(ISBLANK(E41)))'
$cf = $ws.Range("E44").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G44:K44,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G44:K44,"<>" & "")>0,NOT(ISBLANK(E44)))'
$cf = $ws.Range("E45").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G45:K45,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G45:K45,"<>" & "")>0,NOT(ISBLANK(E45)))'
$cf = $ws.Range("E46").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G46:K46,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G46:K46,"<>" & "")>0,NOT(ISBLANK(E46)))'
$cf = $ws.Range("E47").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G47:K47,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G47:K47,"<>" & "")>0,NOT(ISBLANK(E47)))'
$cf = $ws.Range("E48").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G48:K48,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G48:K48,"<>" & "")>0,NOT(ISBLANK(E48)))'
$cf = $ws.Range("E51").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G51:K51,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G51:K51,"<>" & "")>0,NOT(ISBLANK(E51)))'
$cf = $ws.Range("E52").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G52:K52,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G52:K52,"<>" & "")>0,NOT(ISBLANK(E52)))'
$cf = $ws.Range("E53").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G53:K53,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G53:K53,"<>" & "")>0,NOT(ISBLANK(E53)))'
$cf = $ws.Range("E54").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G54:K54,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G54:K54,"<>" & "")>0,NOT(ISBLANK(E54)))'
$cf = $ws.Range("E55").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G55:K55,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G55:K55,"<>" & "")>0,NOT(ISBLANK(E55)))'
$cf = $ws.Range("E58").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G58:K58,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G58:K58,"<>" & "")>0,NOT(ISBLANK(E58)))'
$cf = $ws.Range("E59").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G59:K59,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G59:K59,"<>" & "")>0,NOT(ISBLANK(E59)))'
$cf = $ws.Range("E60").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G60:K60,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G60:K60,"<>" & "")>0,NOT(ISBLANK(E60)))'
$cf = $ws.Range("E61").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G61:K61,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G61:K61,"<>" & "")>0,NOT(ISBLANK(E61)))'
$cf = $ws.Range("E62").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G62:K62,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G62:K62,"<>" & "")>0,NOT(ISBLANK(E62)))'
$cf = $ws.Range("E65").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G65:K65,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G65:K65,"<>" & "")>0,NOT(ISBLANK(E65)))'
$cf = $ws.Range("E66").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G66:K66,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G66:K66,"<>" & "")>0,NOT(ISBLANK(E66)))'
$cf = $ws.Range("E67").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G67:K67,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G67:K67,"<>" & "")>0,NOT(ISBLANK(E67)))'
$cf = $ws.Range("E68").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G68:K68,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G68:K68,"<>" & "")>0,NOT(ISBLANK(E68)))'
$cf = $ws.Range("E69").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G69:K69,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G69:K69,"<>" & "")>0,NOT(ISBLANK(E69)))'
$cf = $ws.Range("E72").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G72:K72,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G72:K72,"<>" & "")>0,NOT(ISBLANK(E72)))'
$cf = $ws.Range("E73").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G73:K73,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G73:K73,"<>" & "")>0,NOT(ISBLANK(E73)))'
$cf = $ws.Range("E74").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G74:K74,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G74:K74,"<>" & "")>0,NOT(ISBLANK(E74)))'
$cf = $ws.Range("E75").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G75:K75,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G75:K75,"<>" & "")>0,NOT(ISBLANK(E75)))'
$cf = $ws.Range("E76").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G76:K76,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G76:K76,"<>" & "")>0,NOT(ISBLANK(E76)))'
$cf = $ws.Range("E79").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G79:K79,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G79:K79,"<>" & "")>0,NOT(ISBLANK(E79)))'
$cf = $ws.Range("E80").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G80:K80,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G80:K80,"<>" & "")>0,NOT(ISBLANK(E80)))'
$cf = $ws.Range("E81").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G81:K81,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G81:K81,"<>" & "")>0,NOT(ISBLANK(E81)))'
$cf = $ws.Range("E82").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G82:K82,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G82:K82,"<>" & "")>0,NOT(ISBLANK(E82)))'
$cf = $ws.Range("E83").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G83:K83,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G83:K83,"<>" & "")>0,NOT(ISBLANK(E83)))'
$cf = $ws.Range("E86").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G86:K86,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G86:K86,"<>" & "")>0,NOT(ISBLANK(E86)))'
$cf = $ws.Range("E87").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G87:K87,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G87:K87,"<>" & "")>0,NOT(ISBLANK(E87)))'
$cf = $ws.Range("E88").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G88:K88,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G88:K88,"<>" & "")>0,NOT(ISBLANK(E88)))'
$cf = $ws.Range("E89").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G89:K89,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G89:K89,"<>" & "")>0,NOT(ISBLANK(E89)))'
$cf = $ws.Range("E90").FormatConditions
$cf.Item(1).Formula1 = '=COUNTIF(G90:K90,"<>" & "")>0'
$cf.Item(2).Formula1 = '=AND(COUNTIF(G90:K90,"<>" & "")>0,NOT(ISBLANK(E90)))'
